$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A5").Value = "Carlos"
$ws.Range("B5").Value = "carlos"
$ws.Range("C5").Value = "2021002252@ifam.edu.br"
$ws.Range("D5").Value = "admin123"

$ws.Hyperlinks.Add($ws.Range("C5"), "mailto:2021002252@ifam.edu.br", "", "", "2021002252@ifam.edu.br")

$ws.Range("A5:D5").Style = $ws.Range("A4:D4").Style
$ws.Range("A5").Font.Name = $ws.Range("A4").Font.Name
